# Update financial figures on slide 1 (observaciones / movimientos refresh).
# These shapes have spAutoFit, so re-writing the run text makes PowerPoint
# recompute the shape height; restore the original EMU height afterwards
# (the magic literal is the point value that round-trips back to exactly
# 307777 EMU given the engine's single-precision float storage).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-RunText($shapeName, $newText) {
    $shp = $s.Shapes.Item($shapeName)
    $shp.TextFrame.TextRange.Runs(1).Text = $newText
    $shp.Height = 24.2344112396
}

Set-RunText "14 CuadroTexto" "24,341"
Set-RunText "16 CuadroTexto" "$ 5,348,224.05"
Set-RunText "26 CuadroTexto" "2,333"
Set-RunText "27 CuadroTexto" "$ 2,492,903.77"
